# Update "想去人数" (number of people interested) counts in both the
# "展览" and "全部类型" worksheets, which contain the same data.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 115
    $ws.Range("F3").Value = 23
    $ws.Range("F4").Value = 974
}
